$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reorder columns B (Client), C (Type), D (Name) -> B (Name), C (Client), D (Type)
# and add a new "Date" column in F, keeping E (File) and G (file format) where they are.

$rows = 2..13
$oldB = @{}
$oldC = @{}
$oldD = @{}
foreach ($r in $rows) {
    $oldB[$r] = $ws.Cells.Item($r, 2).Value2
    $oldC[$r] = $ws.Cells.Item($r, 3).Value2
    $oldD[$r] = $ws.Cells.Item($r, 4).Value2
}

# Headers
$ws.Range("B1").Value2 = "Name"
$ws.Range("C1").Value2 = "Client"
$ws.Range("D1").Value2 = "Type"
$ws.Range("F1").Value2 = "Date"

foreach ($r in $rows) {
    $ws.Cells.Item($r, 2).Value2 = $oldD[$r]
    $ws.Cells.Item($r, 3).Value2 = $oldB[$r]
    $ws.Cells.Item($r, 4).Value2 = $oldC[$r]
    $ws.Cells.Item($r, 6).Value2 = "Feb 29, 2004 (00:00:00 EST)"
}

$ws.Range("H6").Select()
